$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("body")
$ws.Activate()

$ws.Range("C8").Value = '//*[@id="root"]/div/div/div[4]/div/div[1]/img'

$ws.Range("C7").Select()
